$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. "64.659.30") that must
# remain plain text rather than being auto-converted to numbers by Excel.
# Force Text format before assigning, then restore the default style so the
# cell formatting/style index is unchanged from the original file.
$dCells = @('D2','D3','D5','D6','D8','D9','D13','D16','D17','D18','D19','D21','D22','D23','D25','D28','D30','D31','D33','D34','D37','D39','D41','D44','D45','D46','D47')
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '64.659.30'
$ws.Range('D3').Value = '3.442.55'
$ws.Range('D5').Value = '575.66'
$ws.Range('D6').Value = '161.14'
$ws.Range('D8').Value = '3.445.28'
$ws.Range('D9').Value = '0.580'
$ws.Range('D13').Value = '4.035.78'
$ws.Range('D16').Value = '28.25'
$ws.Range('D17').Value = '64.689.36'
$ws.Range('D18').Value = '3.457.64'
$ws.Range('D19').Value = '6.36'
$ws.Range('D21').Value = '385.34'
$ws.Range('D22').Value = '8.17'
$ws.Range('D23').Value = '73.26'
$ws.Range('D25').Value = '0.998'
$ws.Range('D28').Value = '0.181'
$ws.Range('D30').Value = '6.20'
$ws.Range('D31').Value = '1.43'
$ws.Range('D33').Value = '6.59'
$ws.Range('D34').Value = '23.61'
$ws.Range('D37').Value = '163.18'
$ws.Range('D39').Value = '3.006.55'
$ws.Range('D41').Value = '0.0767'
$ws.Range('D44').Value = '42.83'
$ws.Range('D45').Value = '0.0317'
$ws.Range('D46').Value = '0.772'
$ws.Range('D47').Value = '24.81'

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining changed cells (coin name, link, volume%) are plain text already
# and do not require special number-format handling.
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('E6').Value = '  +2.73%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('E9').Value = '  +8.73%  '
$ws.Range('E10').Value = '  -3.20%  '
$ws.Range('E11').Value = '  +1.91%  '
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('E14').Value = '  -2.10%  '
$ws.Range('E15').Value = '  +4.08%  '
$ws.Range('E16').Value = '  +2.81%  '
$ws.Range('E17').Value = '  +1.56%  '
$ws.Range('E18').Value = '  +0.62%  '
$ws.Range('E19').Value = '  -1.31%  '
$ws.Range('E20').Value = '  -0.90%  '
$ws.Range('E21').Value = '  -1.83%  '
$ws.Range('E22').Value = '  -3.91%  '
$ws.Range('E23').Value = '  +1.81%  '
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('E26').Value = '  +13.64%  '
$ws.Range('E27').Value = '  +2.78%  '
$ws.Range('E28').Value = '  +0.29%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  +5.99%  '
$ws.Range('E31').Value = '  +3.71%  '
$ws.Range('E32').Value = '  -0.12%  '
$ws.Range('E33').Value = '  -1.97%  '
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  +3.70%  '
$ws.Range('E37').Value = '  +3.17%  '
$ws.Range('E39').Value = '  +4.89%  '
$ws.Range('E40').Value = '  +1.53%  '
$ws.Range('E41').Value = '  -2.51%  '
$ws.Range('E42').Value = '  -2.99%  '
$ws.Range('E43').Value = '  +3.34%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('E44').Value = '  +2.06%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E45').Value = '  -0.63%  '
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('E47').Value = '  +9.74%  '
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('E49').Value = '  +6.52%  '
$ws.Range('E50').Value = '  +4.14%  '
$ws.Range('E51').Value = '  +3.73%  '
